$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 9499.5
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null
$ws.Range("H107").Value = 849.9231
$ws.Range("I107").Value = 586.6842
$ws.Range("J107").Value = 1564.4286
$ws.Range("K107").Value = 586.6842
$ws.Range("L107").Value = 1564.4286
$ws.Range("M107").Value = 1333.3158
$ws.Range("N107").Value = -5404.4286
$ws.Range("H113").Value = 6833
$ws.Range("I113").Value = 6749.5
$ws.Range("K113").Value = 6749.5
$ws.Range("M113").Value = -3495.5
$ws.Range("H116").Value = 8914.675999999999
$ws.Range("J116").Value = 9936.321
$ws.Range("L116").Value = 9936.321
$ws.Range("N116").Value = -16820.321
$ws.Range("H125").Value = 4241.25
$ws.Range("J125").Value = 4241.25
$ws.Range("L125").Value = 38171.25
$ws.Range("N125").Value = -43091.25
$ws.Range("H137").Value = 2166.0908
$ws.Range("I137").Value = 2100.6667
$ws.Range("K137").Value = 6302.000100000001
$ws.Range("M137").Value = -3752.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3228.9243
$ws.Range("I32").Value = 2377.7932
$ws.Range("J32").Value = 9399.625
$ws.Range("K32").Value = 2377.7932
$ws.Range("L32").Value = 9399.625
$ws.Range("M32").Value = -2090.7932
$ws.Range("N32").Value = -9973.625
$ws.Range("H43").Value = 37866
$ws.Range("J43").Value = 37486.285
$ws.Range("L43").Value = 37486.285
$ws.Range("N43").Value = -38112.285
$ws.Range("H61").Value = 17694314
$ws.Range("I61").Value = 26251760
$ws.Range("K61").Value = 26251760
$ws.Range("M61").Value = -26251548
$ws.Range("H63").Value = 4528.0713
$ws.Range("I63").Value = 4454
$ws.Range("K63").Value = 4454
$ws.Range("M63").Value = -3768
$ws.Range("H66").Value = 4528.0713
$ws.Range("I66").Value = 4454
$ws.Range("K66").Value = 22270
$ws.Range("M66").Value = -18838
$ws.Range("H122").Value = 5665.5
$ws.Range("I122").Value = 5728.3076
$ws.Range("K122").Value = 17184.9228
$ws.Range("M122").Value = -14734.9228
$ws.Range("H132").Value = 2636950
$ws.Range("I132").Value = 5026.4194
$ws.Range("K132").Value = 15079.2582
$ws.Range("M132").Value = -12549.2582
$ws.Range("H136").Value = 17694314
$ws.Range("I136").Value = 26251760
$ws.Range("K136").Value = 78755280
$ws.Range("M136").Value = -78752730

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7002.2666
$ws.Range("I20").Value = 8277.833000000001
$ws.Range("K20").Value = 8277.833000000001
$ws.Range("M20").Value = -8030.833000000001
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35472
$ws.Range("H105").Value = 860870.7
$ws.Range("I105").Value = 1608107.8
$ws.Range("J105").Value = 6885.4287
$ws.Range("K105").Value = 1608107.8
$ws.Range("L105").Value = 6885.4287
$ws.Range("M105").Value = -1606360.8
$ws.Range("N105").Value = -10379.4287
$ws.Range("H134").Value = 5002997.5
$ws.Range("I134").Value = 2564.7334
$ws.Range("J134").Value = 20004296
$ws.Range("K134").Value = 7694.2002
$ws.Range("L134").Value = 60012888
$ws.Range("M134").Value = -5159.2002
$ws.Range("N134").Value = -60017958
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6683066.5
$ws.Range("J16").Value = 2999.5
$ws.Range("L16").Value = 2999.5
$ws.Range("N16").Value = -3573.5
$ws.Range("H58").Value = 2726
$ws.Range("I58").Value = 2726
$ws.Range("K58").Value = 2726
$ws.Range("M58").Value = -2523
$ws.Range("H107").Value = 958.9666999999999
$ws.Range("I107").Value = 658.96
$ws.Range("J107").Value = 2459
$ws.Range("K107").Value = 658.96
$ws.Range("L107").Value = 2459
$ws.Range("M107").Value = 1261.04
$ws.Range("N107").Value = -6299
$ws.Range("H113").Value = 6683066.5
$ws.Range("J113").Value = 2999.5
$ws.Range("L113").Value = 2999.5
$ws.Range("N113").Value = -7339.5
$ws.Range("H125").Value = 78499.25
$ws.Range("J125").Value = 78499.25
$ws.Range("L125").Value = 78499.25
$ws.Range("N125").Value = -83419.25
$ws.Range("H132").Value = 2802.2222
$ws.Range("I132").Value = 2720.524
$ws.Range("J132").Value = 3088.1667
$ws.Range("K132").Value = 8161.572
$ws.Range("L132").Value = 9264.500100000001
$ws.Range("M132").Value = -5631.572
$ws.Range("N132").Value = -14324.5001
$ws.Range("H136").Value = 2726
$ws.Range("I136").Value = 2726
$ws.Range("K136").Value = 8178
$ws.Range("M136").Value = -5628

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 637.875
$ws.Range("J86").Value = 980.6
$ws.Range("L86").Value = 2941.8
$ws.Range("N86").Value = -5313.8
$ws.Range("H89").Value = 637.875
$ws.Range("J89").Value = 980.6
$ws.Range("L89").Value = 8825.4
$ws.Range("N89").Value = -20681.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2772.2666
$ws.Range("I102").Value = 2771.76
$ws.Range("K102").Value = 2771.76
$ws.Range("M102").Value = -1149.76
$ws.Range("H107").Value = 1327.3334
$ws.Range("I107").Value = 1315.3889
$ws.Range("K107").Value = 1315.3889
$ws.Range("M107").Value = 604.6111000000001
$ws.Range("H113").Value = 1854900.2
$ws.Range("I113").Value = 3398.8
$ws.Range("J113").Value = 3706401.5
$ws.Range("K113").Value = 3398.8
$ws.Range("L113").Value = 3706401.5
$ws.Range("M113").Value = -1228.8
$ws.Range("N113").Value = -3710741.5
$ws.Range("H122").Value = 2463.9092
$ws.Range("I122").Value = 2455.889
$ws.Range("K122").Value = 7367.667
$ws.Range("M122").Value = -4917.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17332.834
$ws.Range("I22").Value = 20499.4
$ws.Range("K22").Value = 20499.4
$ws.Range("M22").Value = -20204.4
$ws.Range("H27").Value = 17332.834
$ws.Range("I27").Value = 20499.4
$ws.Range("K27").Value = 20499.4
$ws.Range("M27").Value = -20392.4
$ws.Range("H122").Value = 3642.3333
$ws.Range("I122").Value = 3336.9429
$ws.Range("K122").Value = 10010.8287
$ws.Range("M122").Value = -7560.8287
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 685.0476
$ws.Range("I113").Value = 568.75
$ws.Range("K113").Value = 1706.25
$ws.Range("M113").Value = 463.75
$ws.Range("H122").Value = 2615.15
$ws.Range("I122").Value = 2461.3333
$ws.Range("K122").Value = 7383.999899999999
$ws.Range("M122").Value = -4933.999899999999
$ws.Range("H132").Value = 408243.72
$ws.Range("I132").Value = 7854.9
$ws.Range("K132").Value = 23564.7
$ws.Range("M132").Value = -21034.7
$ws.Range("H136").Value = 1143874.5
$ws.Range("J136").Value = 5012500
$ws.Range("L136").Value = 15037500
$ws.Range("N136").Value = -15042600
